$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp caption (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 09:46"

# --- Update country labels (column A) to reflect new leaderboard ordering ---
$ws.Range("A44").Value = "Filipinas"  # was: India
$ws.Range("A45").Value = "India"  # was: Singapur
$ws.Range("A46").Value = "Singapur"  # was: Filipinas
$ws.Range("A51").Value = "Croacia"  # was: Peru
$ws.Range("A52").Value = "Peru"  # was: Mexico
$ws.Range("A53").Value = "Mexico"  # was: Egipto
$ws.Range("A54").Value = "Egipto"  # was: Barein
$ws.Range("A55").Value = "Barein"  # was: Argentina
$ws.Range("A56").Value = "Argentina"  # was: Hong Kong
$ws.Range("A57").Value = "Hong Kong"  # was: Croacia
$ws.Range("A70").Value = "Letonia"  # was: Bulgaria
$ws.Range("A71").Value = "Bulgaria"  # was: Nueva Zelanda
$ws.Range("A72").Value = "Nueva Zelanda"  # was: Eslovaquia
$ws.Range("A73").Value = "Eslovaquia"  # was: Letonia
$ws.Range("A141").Value = "Nueva Caledonia"  # was: Maldivas
$ws.Range("A142").Value = "Maldivas"  # was: Nueva Caledonia
$ws.Range("A143").Value = "Tanzania"  # was: Etiopia
$ws.Range("A144").Value = "Etiopia"  # was: Tanzania
$ws.Range("A151").Value = "Dominica"  # was: Seychelles
$ws.Range("A153").Value = "Seychelles"  # was: Dominica
$ws.Range("A154").Value = "Niger"  # was: Namibia
$ws.Range("A155").Value = "Namibia"  # was: Bermudas
$ws.Range("A156").Value = "Bermudas"  # was: Benin
$ws.Range("A157").Value = "Benin"  # was: Islas Caimanes
$ws.Range("A158").Value = "Gabon"  # was: Curazao
$ws.Range("A159").Value = "Islas Caimanes"  # was: Gabon
$ws.Range("A160").Value = "Curazao"  # was: Fiyi
$ws.Range("A161").Value = "Fiyi"  # was: Guyana
$ws.Range("A162").Value = "Guyana"  # was: Bahamas
$ws.Range("A163").Value = "Bahamas"  # was: Groenlandia
$ws.Range("A164").Value = "Groenlandia"  # was: Suazilandia
$ws.Range("A165").Value = "Suazilandia"  # was: Guinea
$ws.Range("A166").Value = "Guinea"  # was: Congo
$ws.Range("A167").Value = "Congo"  # was: Santa Sede
$ws.Range("A168").Value = "Santa Sede"  # was: Liberia
$ws.Range("A169").Value = "Liberia"  # was: Mozambique
$ws.Range("A170").Value = "Mozambique"  # was: San Bartolome
$ws.Range("A171").Value = "San Bartolome"  # was: Republica del Chad
$ws.Range("A172").Value = "Republica del Chad"  # was: Republica de Yibuti
$ws.Range("A173").Value = "Republica de Yibuti"  # was: Zambia
$ws.Range("A174").Value = "Zambia"  # was: Birmania
$ws.Range("A175").Value = "Birmania"  # was: Antigua y Barbuda
$ws.Range("A176").Value = "Antigua y Barbuda"  # was: Niger
$ws.Range("A178").Value = "Republica de Africa Central"  # was: Angola
$ws.Range("A179").Value = "Angola"  # was: Republica de Africa Central
$ws.Range("A181").Value = "Cabo Verde"  # was: Gambia
$ws.Range("A183").Value = "Gambia"  # was: Sudan
$ws.Range("A184").Value = "Sudan"  # was: Cabo Verde
$ws.Range("A185").Value = "San Martin (Parte Holandesa)"  # was: Laos
$ws.Range("A186").Value = "Mauritania"  # was: Butan
$ws.Range("A188").Value = "Laos"  # was: Mauritania
$ws.Range("A189").Value = "Butan"  # was: San Martin (Parte Holandesa)
$ws.Range("A190").Value = "Timor Oriental"  # was: Eritrea
$ws.Range("A193").Value = "Belice"  # was: Montserrat
$ws.Range("A194").Value = "San Vicente y las Granadinas"  # was: Belice
$ws.Range("A195").Value = "Montserrat"  # was: San Vicente y las Granadinas
$ws.Range("A200").Value = "Eritrea"  # was: Timor Oriental

# --- Update statistics (columns B-H) with refreshed case counts ---
$ws.Range("B6").Value = 54935
$ws.Range("C6").Value = 54
$ws.Range("E6").Value = 53772

$ws.Range("B8").Value = 33593
$ws.Range("C8").Value = 602
$ws.Range("D8").Value = 3299
$ws.Range("E8").Value = 30130
$ws.Range("G8").Value = 5
$ws.Range("H8").Value = 164

$ws.Range("B15").Value = 5427
$ws.Range("C15").Value = 144
$ws.Range("E15").Value = 5388

$ws.Range("B44").Value = 636
$ws.Range("C44").Value = 84
$ws.Range("D44").Value = 26
$ws.Range("E44").Value = 572
$ws.Range("F44").Value = 1
$ws.Range("G44").Value = 3
$ws.Range("H44").Value = 38

$ws.Range("B45").Value = 562
$ws.Range("C45").Value = 26
$ws.Range("D45").Value = 40
$ws.Range("E45").Value = 512
$ws.Range("F45").Value = 0
$ws.Range("H45").Value = 10

$ws.Range("B46").Value = 558
$ws.Range("D46").Value = 156
$ws.Range("E46").Value = 400
$ws.Range("F46").Value = 17
$ws.Range("H46").Value = 2

$ws.Range("B51").Value = 418
$ws.Range("C51").Value = 36
$ws.Range("D51").Value = 16
$ws.Range("E51").Value = 401
$ws.Range("F51").Value = 6
$ws.Range("H51").Value = 1

$ws.Range("B52").Value = 416
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 1
$ws.Range("E52").Value = 408
$ws.Range("F52").Value = 9
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 7

$ws.Range("B53").Value = 405
$ws.Range("C53").Value = 38
$ws.Range("D53").Value = 4
$ws.Range("E53").Value = 396
$ws.Range("F53").Value = 1
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 5

$ws.Range("B54").Value = 402
$ws.Range("D54").Value = 80
$ws.Range("E54").Value = 302
$ws.Range("F54").Value = 0
$ws.Range("H54").Value = 20

$ws.Range("B55").Value = 392
$ws.Range("D55").Value = 177
$ws.Range("E55").Value = 212
$ws.Range("F55").Value = 2
$ws.Range("H55").Value = 3

$ws.Range("D56").Value = 52
$ws.Range("E56").Value = 329
$ws.Range("F56").Value = 0
$ws.Range("H56").Value = 6

$ws.Range("B57").Value = 387
$ws.Range("D57").Value = 102
$ws.Range("E57").Value = 281
$ws.Range("F57").Value = 4
$ws.Range("H57").Value = 4

$ws.Range("B70").Value = 221
$ws.Range("C70").Value = 24
$ws.Range("D70").Value = 1
$ws.Range("E70").Value = 220
$ws.Range("F70").Value = 0
$ws.Range("H70").Value = 0

$ws.Range("B71").Value = 220
$ws.Range("C71").Value = 2
$ws.Range("D71").Value = 4
$ws.Range("E71").Value = 213
$ws.Range("F71").Value = 8
$ws.Range("H71").Value = 3

$ws.Range("B72").Value = 205
$ws.Range("D72").Value = 22
$ws.Range("E72").Value = 183
$ws.Range("F72").Value = 0

$ws.Range("B73").Value = 204
$ws.Range("D73").Value = 7
$ws.Range("E73").Value = 197
$ws.Range("F73").Value = 2

$ws.Range("B74").Value = 195
$ws.Range("C74").Value = 4
$ws.Range("D74").Value = 43
$ws.Range("F74").Value = 6

$ws.Range("B100").Value = 86
$ws.Range("C100").Value = 5
$ws.Range("E100").Value = 64
$ws.Range("F100").Value = 2

$ws.Range("B141").Value = 14
$ws.Range("C141").Value = 4
$ws.Range("D141").Value = 0
$ws.Range("E141").Value = 14

$ws.Range("B142").Value = 13
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 5
$ws.Range("E142").Value = 8

$ws.Range("C154").Value = 4
$ws.Range("D154").Value = 0
$ws.Range("E154").Value = 6
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 1

$ws.Range("B155").Value = 7
$ws.Range("D155").Value = 2
$ws.Range("E155").Value = 5

$ws.Range("E157").Value = 6
$ws.Range("H157").Value = 0

$ws.Range("B160").Value = 6
$ws.Range("C160").Value = 0
$ws.Range("H160").Value = 1

$ws.Range("C161").Value = 1
$ws.Range("E161").Value = 5
$ws.Range("H161").Value = 0

$ws.Range("D162").Value = 0
$ws.Range("H162").Value = 1

$ws.Range("D163").Value = 1
$ws.Range("E163").Value = 4

$ws.Range("B164").Value = 5
$ws.Range("D164").Value = 2
$ws.Range("E164").Value = 3

$ws.Range("B168").Value = 4
$ws.Range("E168").Value = 4

